$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 0.06414394825696945
$ws.Cells.Item(2, 2).Value = 0.9805817008018494
$ws.Cells.Item(2, 3).Value = 0.003310504369437695
$ws.Cells.Item(2, 4).Value = 0.9994794726371765

$ws.Cells.Item(3, 1).Value = 0.01024373713880777
$ws.Cells.Item(3, 2).Value = 0.9979604482650757
$ws.Cells.Item(3, 3).Value = 0.002520929789170623
$ws.Cells.Item(3, 4).Value = 0.9994794726371765

$ws.Cells.Item(4, 1).Value = 0.006017731502652168
$ws.Cells.Item(4, 2).Value = 0.9985978007316589
$ws.Cells.Item(4, 3).Value = 0.0008317870087921619
$ws.Cells.Item(4, 4).Value = 0.9998264908790588

$ws.Cells.Item(5, 1).Value = 0.003417884465306997
$ws.Cells.Item(5, 2).Value = 0.9991501569747925
$ws.Cells.Item(5, 3).Value = 0.0006752578774467111
$ws.Cells.Item(5, 4).Value = 0.999884307384491

$ws.Cells.Item(6, 1).Value = 0.002150098793208599
$ws.Cells.Item(6, 2).Value = 0.9995113611221313
$ws.Cells.Item(6, 3).Value = 0.0001991159369936213
$ws.Cells.Item(6, 4).Value = 0.9998264908790588

$ws.Cells.Item(7, 1).Value = 0.001314928871579468
$ws.Cells.Item(7, 2).Value = 0.9996600747108459
$ws.Cells.Item(7, 3).Value = 0.00005622482058242895
$ws.Cells.Item(7, 4).Value = 1

$ws.Cells.Item(8, 1).Value = 0.001312813838012516
$ws.Cells.Item(8, 2).Value = 0.9997025728225708
$ws.Cells.Item(8, 3).Value = 0.00007380757597275078
$ws.Cells.Item(8, 4).Value = 1

$ws.Cells.Item(9, 1).Value = 0.001131516881287098
$ws.Cells.Item(9, 2).Value = 0.9997662901878357
$ws.Cells.Item(9, 3).Value = 0.0002006521390285343
$ws.Cells.Item(9, 4).Value = 0.9998264908790588

$ws.Cells.Item(10, 1).Value = 0.001083581824786961
$ws.Cells.Item(10, 2).Value = 0.9997237920761108
$ws.Cells.Item(10, 3).Value = 0.00002121128090948332
$ws.Cells.Item(10, 4).Value = 1

$ws.Cells.Item(11, 1).Value = 0.0007308748317882419
$ws.Cells.Item(11, 2).Value = 0.9997662901878357
$ws.Cells.Item(11, 3).Value = 0.000001265031983166409
$ws.Cells.Item(11, 4).Value = 1

$ws.Cells.Item(12, 1).Value = 0.001217738958075643
$ws.Cells.Item(12, 2).Value = 0.9997875690460205
$ws.Cells.Item(12, 3).Value = 0.0000004983010057912907
$ws.Cells.Item(12, 4).Value = 1

$ws.Cells.Item(13, 1).Value = 0.000476664979942143
$ws.Cells.Item(13, 2).Value = 0.9998300075531006
$ws.Cells.Item(13, 3).Value = 0.0000005235870048636571
$ws.Cells.Item(13, 4).Value = 1

$ws.Cells.Item(14, 1).Value = 0.00007042202196316794
$ws.Cells.Item(14, 2).Value = 0.99997878074646
$ws.Cells.Item(14, 3).Value = 0.0000001075793250038259
$ws.Cells.Item(14, 4).Value = 1

$ws.Cells.Item(15, 1).Value = 0.0007841288461349905
$ws.Cells.Item(15, 2).Value = 0.9998087882995605
$ws.Cells.Item(15, 3).Value = 0.00000004218784965814848
$ws.Cells.Item(15, 4).Value = 1

$ws.Cells.Item(16, 1).Value = 0.001342142815701663
$ws.Cells.Item(16, 2).Value = 0.9997662901878357
$ws.Cells.Item(16, 3).Value = 0.00001468345180910546
$ws.Cells.Item(16, 4).Value = 1

$ws.Cells.Item(17, 1).Value = 0.000276384613243863
$ws.Cells.Item(17, 2).Value = 0.9999362826347351
$ws.Cells.Item(17, 3).Value = 0.000001209832021231705
$ws.Cells.Item(17, 4).Value = 1

$ws.Cells.Item(18, 1).Value = 0.0005018216324970126
$ws.Cells.Item(18, 2).Value = 0.9998087882995605
$ws.Cells.Item(18, 3).Value = 0.000001119099465540785
$ws.Cells.Item(18, 4).Value = 1

$ws.Cells.Item(19, 1).Value = 0.0003714240156114101
$ws.Cells.Item(19, 2).Value = 0.9998512864112854
$ws.Cells.Item(19, 3).Value = 0.00000001598033882999061
$ws.Cells.Item(19, 4).Value = 1

$ws.Cells.Item(20, 1).Value = 0.0006294162012636662
$ws.Cells.Item(20, 2).Value = 0.9998512864112854
$ws.Cells.Item(20, 3).Value = 0.00000001633759261210344
$ws.Cells.Item(20, 4).Value = 1

$ws.Cells.Item(21, 1).Value = 0.0001296165864914656
$ws.Cells.Item(21, 2).Value = 0.9999362826347351
$ws.Cells.Item(21, 3).Value = 0.00001486019846197451
$ws.Cells.Item(21, 4).Value = 1

$ws.Cells.Item(22, 1).Value = 0.0003832871443592012
$ws.Cells.Item(22, 2).Value = 0.9999150037765503
$ws.Cells.Item(22, 3).Value = 0.000001090185719476722
$ws.Cells.Item(22, 4).Value = 1

$ws.Cells.Item(23, 1).Value = 0.000565042719244957
$ws.Cells.Item(23, 2).Value = 0.9998512864112854
$ws.Cells.Item(23, 3).Value = 0.0000007335787586271181
$ws.Cells.Item(23, 4).Value = 1

$ws.Cells.Item(24, 1).Value = 0.0001239280245499685
$ws.Cells.Item(24, 2).Value = 0.9999575018882751
$ws.Cells.Item(24, 3).Value = 0.00000005421972204544545
$ws.Cells.Item(24, 4).Value = 1

$ws.Cells.Item(25, 1).Value = 0.0006511447136290371
$ws.Cells.Item(25, 2).Value = 0.9998087882995605
$ws.Cells.Item(25, 3).Value = 0.000000003288549654811845
$ws.Cells.Item(25, 4).Value = 1

$ws.Cells.Item(26, 1).Value = 0.0002296371967531741
$ws.Cells.Item(26, 2).Value = 0.9999150037765503
$ws.Cells.Item(26, 3).Value = 0.000000008383232774633598
$ws.Cells.Item(26, 4).Value = 1

$ws.Cells.Item(27, 1).Value = 0.0001516987103968859
$ws.Cells.Item(27, 2).Value = 0.9999362826347351
$ws.Cells.Item(27, 3).Value = 0.000000002206162363904696
$ws.Cells.Item(27, 4).Value = 1

$ws.Cells.Item(28, 1).Value = 0.0003954574931412935
$ws.Cells.Item(28, 2).Value = 0.9999362826347351
$ws.Cells.Item(28, 3).Value = 0.000000000551543255511433
$ws.Cells.Item(28, 4).Value = 1

$ws.Cells.Item(29, 1).Value = 0.00006237801426323131
$ws.Cells.Item(29, 2).Value = 0.99997878074646
$ws.Cells.Item(29, 3).Value = 0.0000000001172030250629064
$ws.Cells.Item(29, 4).Value = 1

$ws.Cells.Item(30, 1).Value = 0.00009686122211860493
$ws.Cells.Item(30, 2).Value = 0.9999362826347351
$ws.Cells.Item(30, 3).Value = 0.00000000002757718532697684
$ws.Cells.Item(30, 4).Value = 1

$ws.Cells.Item(31, 1).Value = 0.00007637272938154638
$ws.Cells.Item(31, 2).Value = 0.99997878074646
$ws.Cells.Item(31, 3).Value = 0.0000000002137230531440437
$ws.Cells.Item(31, 4).Value = 1

$ws.Cells.Item(32, 1).Value = 0.0002198486181441694
$ws.Cells.Item(32, 2).Value = 0.9999362826347351
$ws.Cells.Item(32, 3).Value = 0.0006672106101177633
$ws.Cells.Item(32, 4).Value = 0.999884307384491

$ws.Cells.Item(33, 1).Value = 0.00002413172660453711
$ws.Cells.Item(33, 2).Value = 0.99997878074646
$ws.Cells.Item(33, 3).Value = 0.00000001406351124444427
$ws.Cells.Item(33, 4).Value = 1

$ws.Cells.Item(34, 1).Value = 0.0004056181060150266
$ws.Cells.Item(34, 2).Value = 0.9999362826347351
$ws.Cells.Item(34, 3).Value = 0.000000001068611199528391
$ws.Cells.Item(34, 4).Value = 1

$ws.Cells.Item(35, 1).Value = 0.000229976067203097
$ws.Cells.Item(35, 2).Value = 0.9998937845230103
$ws.Cells.Item(35, 3).Value = 0.000006698771358060185
$ws.Cells.Item(35, 4).Value = 1

$ws.Cells.Item(36, 1).Value = 0.0003630123101174831
$ws.Cells.Item(36, 2).Value = 0.9999362826347351
$ws.Cells.Item(36, 3).Value = 0.000000004494995931025869
$ws.Cells.Item(36, 4).Value = 1

$ws.Cells.Item(37, 1).Value = 0.0006045158370397985
$ws.Cells.Item(37, 2).Value = 0.9998725056648254
$ws.Cells.Item(37, 3).Value = 0.000000004488059257568011
$ws.Cells.Item(37, 4).Value = 1

$ws.Cells.Item(38, 1).Value = 0.0001424906076863408
$ws.Cells.Item(38, 2).Value = 0.9999575018882751
$ws.Cells.Item(38, 3).Value = 0.000000001137556604469125
$ws.Cells.Item(38, 4).Value = 1

$ws.Cells.Item(39, 1).Value = 0.000157424874487333
$ws.Cells.Item(39, 2).Value = 0.99997878074646
$ws.Cells.Item(39, 3).Value = 0.000000009610296558548725
$ws.Cells.Item(39, 4).Value = 1

$ws.Cells.Item(40, 1).Value = 0.0000421263393945992
$ws.Cells.Item(40, 2).Value = 0.99997878074646
$ws.Cells.Item(40, 3).Value = 0.000001172846623376245
$ws.Cells.Item(40, 4).Value = 1

$ws.Cells.Item(41, 1).Value = 0.00001816760232031811
$ws.Cells.Item(41, 2).Value = 1
$ws.Cells.Item(41, 3).Value = 0.0000000008824637132143209
$ws.Cells.Item(41, 4).Value = 1

$ws.Cells.Item(42, 1).Value = 0.0002367770939599723
$ws.Cells.Item(42, 2).Value = 0.9999362826347351
$ws.Cells.Item(42, 3).Value = 0.0000000001103086996745439
$ws.Cells.Item(42, 4).Value = 1

$ws.Cells.Item(43, 1).Value = 0.0001721723383525386
$ws.Cells.Item(43, 2).Value = 0.9999575018882751
$ws.Cells.Item(43, 3).Value = 0.00000001402824434393324
$ws.Cells.Item(43, 4).Value = 1

$ws.Cells.Item(44, 1).Value = 0.0002920124970842153
$ws.Cells.Item(44, 2).Value = 0.9999362826347351
$ws.Cells.Item(44, 3).Value = 0.00000001189156240144484
$ws.Cells.Item(44, 4).Value = 1

$ws.Cells.Item(45, 1).Value = 0.0003606589161790907
$ws.Cells.Item(45, 2).Value = 0.9998937845230103
$ws.Cells.Item(45, 3).Value = 0.0000000008135228046768361
$ws.Cells.Item(45, 4).Value = 1

$ws.Cells.Item(46, 1).Value = 0.0003097125445492566
$ws.Cells.Item(46, 2).Value = 0.9999362826347351
$ws.Cells.Item(46, 3).Value = 0.0000000002275116484096174
$ws.Cells.Item(46, 4).Value = 1

$ws.Cells.Item(47, 1).Value = 0.000005278584922052687
$ws.Cells.Item(47, 2).Value = 1
$ws.Cells.Item(47, 3).Value = 0.00000000006204865787839964
$ws.Cells.Item(47, 4).Value = 1

$ws.Cells.Item(48, 1).Value = 0.001208834000863135
$ws.Cells.Item(48, 2).Value = 0.9997875690460205
$ws.Cells.Item(48, 3).Value = 0.000002724506430240581
$ws.Cells.Item(48, 4).Value = 1

$ws.Cells.Item(49, 1).Value = 0.0003222219238523394
$ws.Cells.Item(49, 2).Value = 0.9999362826347351
$ws.Cells.Item(49, 3).Value = 0.00000009007789714132741
$ws.Cells.Item(49, 4).Value = 1

$ws.Cells.Item(50, 1).Value = 0.00005049924584454857
$ws.Cells.Item(50, 2).Value = 0.99997878074646
$ws.Cells.Item(50, 3).Value = 0.000000001247862035924641
$ws.Cells.Item(50, 4).Value = 1

$ws.Cells.Item(51, 1).Value = 0.0001154916317318566
$ws.Cells.Item(51, 2).Value = 0.99997878074646
$ws.Cells.Item(51, 3).Value = 0.000000001247862035924641
$ws.Cells.Item(51, 4).Value = 1

